# Auto-generated Excel COM-interop script
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H-N)
# across multiple Leve-profit worksheets to reflect refreshed market data.

$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 3557.1667
$ws.Range("I15").Value = 3557.1667
$ws.Range("K15").Value = 10671.5001
$ws.Range("M15").Value = -10502.5001
$ws.Range("H29").Value = 161.25
$ws.Range("I29").Value = 161.25
$ws.Range("K29").Value = 483.75
$ws.Range("M29").Value = -202.75
$ws.Range("H41").Value = 130.2
$ws.Range("I41").Value = 133
$ws.Range("J41").Value = 126
$ws.Range("K41").Value = 133
$ws.Range("L41").Value = 126
$ws.Range("M41").Value = 307
$ws.Range("N41").Value = -1006
$ws.Range("H62").Value = 8228.666999999999
$ws.Range("I62").Value = 6124.75
$ws.Range("K62").Value = 6124.75
$ws.Range("M62").Value = -5500.75
$ws.Range("H64").Value = 9799.799999999999
$ws.Range("I64").Value = 9666.333000000001
$ws.Range("K64").Value = 9666.333000000001
$ws.Range("M64").Value = -9418.333000000001
$ws.Range("H65").Value = 8228.666999999999
$ws.Range("I65").Value = 6124.75
$ws.Range("K65").Value = 30623.75
$ws.Range("M65").Value = -27503.75
$ws.Range("H67").Value = 9799.799999999999
$ws.Range("I67").Value = 9666.333000000001
$ws.Range("K67").Value = 9666.333000000001
$ws.Range("M67").Value = -8808.333000000001
$ws.Range("H86").Value = 1929.8572
$ws.Range("I86").Value = 1558.2222
$ws.Range("J86").Value = 2598.8
$ws.Range("K86").Value = 1558.2222
$ws.Range("L86").Value = 2598.8
$ws.Range("N86").Value = -4844.8
$ws.Range("M86").Value = -435.2221999999999
$ws.Range("H89").Value = 1929.8572
$ws.Range("I89").Value = 1558.2222
$ws.Range("J89").Value = 2598.8
$ws.Range("K89").Value = 7791.111
$ws.Range("L89").Value = 12994
$ws.Range("N89").Value = -24226
$ws.Range("M89").Value = -2175.111
$ws.Range("H98").Value = 2639.5715
$ws.Range("I98").Value = 2580.5
$ws.Range("K98").Value = 2580.5
$ws.Range("M98").Value = -1082.5
$ws.Range("H122").Value = 2639.5715
$ws.Range("I122").Value = 2580.5
$ws.Range("K122").Value = 7741.5
$ws.Range("M122").Value = -5291.5
$ws.Range("H132").Value = 11323.24
$ws.Range("I132").Value = 11575.2
$ws.Range("K132").Value = 34725.60000000001
$ws.Range("M132").Value = -32195.60000000001
$ws.Range("H137").Value = 1250
$ws.Range("I137").Value = 1250
$ws.Range("K137").Value = 3750
$ws.Range("M137").Value = -1200
$ws.Range("H138").Value = 5788
$ws.Range("J138").Value = 7454.5454
$ws.Range("L138").Value = 22363.6362
$ws.Range("N138").Value = -32643.6362

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 152.5
$ws.Range("J3").Value = 105
$ws.Range("L3").Value = 105
$ws.Range("N3").Value = -335
$ws.Range("H4").Value = 177.90909
$ws.Range("I4").Value = 203.5
$ws.Range("K4").Value = 203.5
$ws.Range("M4").Value = -87.5
$ws.Range("H32").Value = 7722844.5
$ws.Range("I32").Value = 25796.2
$ws.Range("J32").Value = 12533500
$ws.Range("K32").Value = 25796.2
$ws.Range("L32").Value = 12533500
$ws.Range("M32").Value = -25509.2
$ws.Range("N32").Value = -12534074
$ws.Range("H33").Value = 4026
$ws.Range("I33").Value = 4026
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 4026
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = -3697
$ws.Range("H44").Value = 7171
$ws.Range("J44").Value = 7171
$ws.Range("L44").Value = 7171
$ws.Range("N44").Value = -8147
$ws.Range("H45").Value = 2899.8
$ws.Range("I45").Value = 1833
$ws.Range("J45").Value = 4500
$ws.Range("K45").Value = 1833
$ws.Range("L45").Value = 4500
$ws.Range("M45").Value = -1456
$ws.Range("N45").Value = -5254
$ws.Range("H55").Value = 14280
$ws.Range("J55").Value = 20000
$ws.Range("L55").Value = 20000
$ws.Range("N55").Value = -20630
$ws.Range("H63").Value = 3462.5557
$ws.Range("I63").Value = 2645.375
$ws.Range("K63").Value = 2645.375
$ws.Range("M63").Value = -1959.375
$ws.Range("H66").Value = 3462.5557
$ws.Range("I66").Value = 2645.375
$ws.Range("K66").Value = 13226.875
$ws.Range("M66").Value = -9794.875
$ws.Range("H74").Value = 4249.375
$ws.Range("I74").Value = 3285.1428
$ws.Range("K74").Value = 3285.1428
$ws.Range("M74").Value = -2411.1428
$ws.Range("H77").Value = 4249.375
$ws.Range("I77").Value = 3285.1428
$ws.Range("K77").Value = 16425.714
$ws.Range("M77").Value = -12057.714
$ws.Range("H102").Value = 3978.75
$ws.Range("I102").Value = 2473.5
$ws.Range("K102").Value = 2473.5
$ws.Range("M102").Value = -851.5
$ws.Range("H122").Value = 1396.1428
$ws.Range("J122").Value = 1499.5
$ws.Range("L122").Value = 4498.5
$ws.Range("N122").Value = -9398.5
$ws.Range("H132").Value = 897.9091
$ws.Range("I132").Value = 897.9091
$ws.Range("K132").Value = 2693.7273
$ws.Range("M132").Value = -163.7273
$ws.Range("N33").ClearContents()

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("H35").Value = 68599.2
$ws.Range("J35").Value = 81999.25
$ws.Range("L35").Value = 81999.25
$ws.Range("N35").Value = -82619.25
$ws.Range("H105").Value = 1500
$ws.Range("I105").Value = 1500
$ws.Range("K105").Value = 1500
$ws.Range("M105").Value = 247
$ws.Range("M25").ClearContents()

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2994.375
$ws.Range("I16").Value = 2616.5
$ws.Range("K16").Value = 2616.5
$ws.Range("M16").Value = -2329.5
$ws.Range("H22").Value = 4534.7144
$ws.Range("I22").Value = 4179.2
$ws.Range("K22").Value = 4179.2
$ws.Range("M22").Value = -3829.2
$ws.Range("H94").Value = 5742.2856
$ws.Range("I94").Value = 2530.6667
$ws.Range("K94").Value = 2530.6667
$ws.Range("M94").Value = -2079.6667
$ws.Range("H99").Value = 5113.231
$ws.Range("I99").Value = 4225.727
$ws.Range("K99").Value = 4225.727
$ws.Range("M99").Value = -2727.727
$ws.Range("H113").Value = 2994.375
$ws.Range("I113").Value = 2616.5
$ws.Range("K113").Value = 2616.5
$ws.Range("M113").Value = -446.5
$ws.Range("H126").Value = 5113.231
$ws.Range("I126").Value = 4225.727
$ws.Range("K126").Value = 12677.181
$ws.Range("M126").Value = -10207.181
$ws.Range("H138").Value = 47000
$ws.Range("I138").Value = 47000
$ws.Range("K138").Value = 47000
$ws.Range("M138").Value = -41860

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 127742.875
$ws.Range("I11").Value = 145950.42
$ws.Range("K11").Value = 437851.26
$ws.Range("M11").Value = -437711.26
$ws.Range("H131").Value = 2000
$ws.Range("J131").Value = 1800
$ws.Range("L131").Value = 5400
$ws.Range("N131").Value = -15480

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 94.25
$ws.Range("I2").Value = 42.6
$ws.Range("J2").Value = 249.2
$ws.Range("K2").Value = 42.6
$ws.Range("L2").Value = 249.2
$ws.Range("M2").Value = 70.40000000000001
$ws.Range("N2").Value = -475.2
$ws.Range("H9").Value = 914.8182
$ws.Range("I9").Value = 206.3
$ws.Range("J9").Value = 8000
$ws.Range("K9").Value = 206.3
$ws.Range("L9").Value = 8000
$ws.Range("M9").Value = -36.30000000000001
$ws.Range("N9").Value = -8340

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 4243.2666
$ws.Range("I46").Value = 1274.8334
$ws.Range("J46").Value = 6222.222
$ws.Range("K46").Value = 1274.8334
$ws.Range("L46").Value = 6222.222
$ws.Range("M46").Value = -1086.8334
$ws.Range("N46").Value = -6598.222
$ws.Range("H100").Value = 5656.684
$ws.Range("I100").Value = 2952.4546
$ws.Range("K100").Value = 2952.4546
$ws.Range("M100").Value = -2411.4546
$ws.Range("H122").Value = 4868.154

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 6166.6665
$ws.Range("I2").Value = 2500
$ws.Range("K2").Value = 2500
$ws.Range("M2").Value = -2388
$ws.Range("H4").Value = 413863.5
$ws.Range("I4").Value = 550068
$ws.Range("K4").Value = 550068
$ws.Range("M4").Value = -549955
$ws.Range("H62").Value = 9540.571
$ws.Range("I62").Value = 6261.3335
$ws.Range("K62").Value = 6261.3335
$ws.Range("M62").Value = -5637.3335
$ws.Range("H65").Value = 9540.571
$ws.Range("I65").Value = 6261.3335
$ws.Range("K65").Value = 31306.6675
$ws.Range("M65").Value = -28186.6675
$ws.Range("H126").Value = 3444.037
$ws.Range("I126").Value = 2924.875
$ws.Range("J126").Value = 4199.1816
$ws.Range("K126").Value = 8774.625
$ws.Range("L126").Value = 12597.5448
$ws.Range("M126").Value = -6304.625
$ws.Range("N126").Value = -17537.5448
$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("H132").Value = 2790.818
$ws.Range("I132").Value = 2790.818
$ws.Range("K132").Value = 8372.454000000002
$ws.Range("M132").Value = -5842.454000000002
$ws.Range("H136").Value = 5566.7393
$ws.Range("I136").Value = 4902.55
$ws.Range("J136").Value = 9994.666999999999
$ws.Range("K136").Value = 14707.65
$ws.Range("L136").Value = 29984.001
$ws.Range("M136").Value = -12157.65
$ws.Range("N136").Value = -35084.001
$ws.Range("N131").ClearContents()

"done"